# Update "previsao_retorno" worksheet (Resumo_por_Cliente) with refreshed
# data for Bibi PF customers: recomputed "meses sem comprar" labels and
# updated probability / purchase-count / date-window figures for a few
# clients whose records moved since the last extraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose only change is the "situacao" label (col J) shifting by
#     one tenth of a month because the report was regenerated one day later.
$ws.Range("J4").Value   = "INATIVO - 37.8 meses sem comprar"
$ws.Range("J9").Value   = "INATIVO - 20.8 meses sem comprar"
$ws.Range("J11").Value  = "INATIVO - 6.7 meses sem comprar"
$ws.Range("J12").Value  = "INATIVO - 1.8 meses sem comprar"
$ws.Range("J17").Value  = "INATIVO - 4.2 meses sem comprar"
$ws.Range("J24").Value  = "INATIVO - 30.9 meses sem comprar"
$ws.Range("J34").Value  = "INATIVO - 9.3 meses sem comprar"
$ws.Range("J37").Value  = "INATIVO - 28.5 meses sem comprar"
$ws.Range("J38").Value  = "INATIVO - 8.2 meses sem comprar"
$ws.Range("J42").Value  = "INATIVO - 36.1 meses sem comprar"
$ws.Range("J46").Value  = "INATIVO - 9.2 meses sem comprar"
$ws.Range("J47").Value  = "INATIVO - 27.5 meses sem comprar"
$ws.Range("J52").Value  = "INATIVO - 8.5 meses sem comprar"
$ws.Range("J54").Value  = "INATIVO - 11.7 meses sem comprar"
$ws.Range("J73").Value  = "INATIVO - 14.3 meses sem comprar"
$ws.Range("J74").Value  = "INATIVO - 13.1 meses sem comprar"
$ws.Range("J79").Value  = "INATIVO - 23.6 meses sem comprar"
$ws.Range("J87").Value  = "INATIVO - 24.0 meses sem comprar"
$ws.Range("J91").Value  = "INATIVO - 23.5 meses sem comprar"
$ws.Range("J95").Value  = "INATIVO - 13.5 meses sem comprar"
$ws.Range("J97").Value  = "INATIVO - 12.9 meses sem comprar"
$ws.Range("J107").Value = "INATIVO - 24.6 meses sem comprar"
$ws.Range("J111").Value = "INATIVO - 0.7 meses sem comprar"
$ws.Range("J118").Value = "INATIVO - 23.1 meses sem comprar"
$ws.Range("J120").Value = "INATIVO - 8.8 meses sem comprar"

# --- id_cliente 5322 (row 51): new purchase pushed prob/date window forward.
$ws.Range("B51").Value = 0.25
$ws.Range("D51").Value = 0.5
$ws.Range("E51").Value = 18
$ws.Range("F51").Value = 0.5
$ws.Range("H51").Value = 45887.41381944445
$ws.Range("I51").Value = 45979.41381944445

# --- id_cliente 5988 (row 60): regularidade/padrao recalculated.
$ws.Range("C60").Value = 0.33
$ws.Range("D60").Value = 0.33
$ws.Range("E60").Value = 35
$ws.Range("F60").Value = 0.33
$ws.Range("H60").Value = 45887.61050925926
$ws.Range("I60").Value = 45902.61050925926

# --- id_cliente 28458 (row 123): BIBI CELL itself, now fully ATIVO (100%).
$ws.Range("B123").Value = 1
$ws.Range("C123").Value = 1
$ws.Range("E123").Value = 18115
$ws.Range("H123").Value = 45887.79111111111
$ws.Range("I123").Value = 45888.79111111111
